$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("browser" row) switches driver from chrome to firefox
$ws.Range("F2").Value = "firefox"

# The "path" column (G) - header and its value - is removed from the sheet
$ws.Range("G1").ClearContents()
$ws.Range("G2").ClearContents()

# Update the active selection to F4
$ws.Range("F4").Select()
